$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text + column widths ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("I2").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.md"
$zh.Range("I3").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.md"
$zh.Range("J2").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.zh-cn.xlf"
$zh.Range("J3").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 21:11:32"
$zh.Range("K3").Value = "2016-08-30 21:11:32"
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("I2").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.md"
$de.Range("I3").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.md"
$de.Range("J2").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.de-de.xlf"
$de.Range("J3").Value = "a64fdfab-c11e-4542-aa3d-e5683b07d293.4c04c08c4c27ec13b3c2344d1a97a843ae0247b4.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 21:11:40"
$de.Range("K3").Value = "2016-08-30 21:11:40"
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
